$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 3) imported from a sample grave/martyr record
$ws.Range("A3").Value = "Đỗ Văn Cừ "
$ws.Range("B3").Value = 1928
$ws.Range("C3").Value = "Trung đội phó"
$ws.Range("D3").Value = "Trung đội trưởng`n"
$ws.Range("D3").WrapText = $true
$ws.Range("E3").Value = "18-6-1969"
$ws.Range("F3").Value = "ZB 164b"
$ws.Range("G3").Value = "135TTga/04-04-1957"
$ws.Range("H3").Value = "Anh ruột"
$ws.Range("I3").Value = "Lê Hồng Quân"
$ws.Range("L3").Value = "1762188222-3f8776.png"

# Columns widen slightly to accommodate the new content
$ws.Columns("I:I").ColumnWidth = 11.4
$ws.Columns("L:L").ColumnWidth = 21.1
$ws.Columns("M:M").ColumnWidth = 7.7

# Move active selection (matches where the user clicked next)
$ws.Range("K14").Select()
